$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.2946457862854
$ws.Range("B1").Value = 2.427496194839478
$ws.Range("C1").Value = 1.772054195404053
$ws.Range("D1").Value = 1.641387581825256
$ws.Range("E1").Value = 1.555007576942444
